$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55: fill in SPY OPEN/HIGH/LOW/CLOSE/VOLUME/SDs values
$ws.Range("C55").Value = 585.58000000000004
$ws.Range("D55").Value = 590.97
$ws.Range("E55").Value = 585.1
$ws.Range("F55").Value = 590.46
$ws.Range("G55").Value = 71119088
$ws.Range("H55").Value = 0.51248766194008022

# Row 56: new EOD/0DTE dates plus the full options-chain metrics snapshot
$ws.Range("A56").Value = 45792
$ws.Range("B56").Value = 45793
$ws.Range("I56").Value = 0.1782
$ws.Range("J56").Value = 3.7
$ws.Range("K56").Value = 600
$ws.Range("L56").Value = 353871600
$ws.Range("M56").Value = -11028
$ws.Range("N56").Value = 1097
$ws.Range("O56").Value = 12125
$ws.Range("P56").Value = 0.11355252885958711
$ws.Range("Q56").Value = 0.18414972676038049
$ws.Range("R56").Value = 45793
$ws.Range("S56").Value = 0.18414972676038049
$ws.Range("T56").Value = 45800
$ws.Range("U56").Value = 0.071891077867353284
$ws.Range("V56").Value = 45828
$ws.Range("W56").Value = 0.34061602826300047
$ws.Range("X56").Value = 15
$ws.Range("Y56").Value = 590
$ws.Range("Z56").Value = 283450160
$ws.Range("AA56").Value = -12757
$ws.Range("AB56").Value = 23172
$ws.Range("AC56").Value = 35929
$ws.Range("AD56").Value = 0.09095525742572895
$ws.Range("AE56").Value = 0.53833815101338367
$ws.Range("AF56").Value = 45793
$ws.Range("AG56").Value = 0.53833815101338367
$ws.Range("AH56").Value = 45807
$ws.Range("AI56").Value = 0.071388752603142147
$ws.Range("AJ56").Value = 45828
$ws.Range("AK56").Value = 0.10913830108804126
$ws.Range("AL56").Value = 17.333333333333332
$ws.Range("AM56").Value = 595
$ws.Range("AN56").Value = 270300765
$ws.Range("AO56").Value = 24847
$ws.Range("AP56").Value = 1187
$ws.Range("AQ56").Value = 26034
$ws.Range("AR56").Value = 0.086735797443001847
$ws.Range("AS56").Value = 0.37098076437125194
$ws.Range("AT56").Value = 45793
$ws.Range("AU56").Value = 0.37098076437125194
$ws.Range("AV56").Value = 45807
$ws.Range("AW56").Value = 0.19810047395701721
$ws.Range("AX56").Value = 45828
$ws.Range("AY56").Value = 0.097663496681068218
$ws.Range("AZ56").Value = 17.333333333333332
$ws.Range("BA56").Value = 605
$ws.Range("BB56").Value = 182293760
$ws.Range("BC56").Value = 25160
$ws.Range("BD56").Value = 203
$ws.Range("BE56").Value = 25363
$ws.Range("BF56").Value = 0.058495560093894637
$ws.Range("BG56").Value = 0.037494815172507724
$ws.Range("BH56").Value = 45798
$ws.Range("BI56").Value = 0.058082449591095098
$ws.Range("BJ56").Value = 45807
$ws.Range("BK56").Value = 0.24081124751592592
$ws.Range("BL56").Value = 45828
$ws.Range("BM56").Value = 0.35115675822917247
$ws.Range("BN56").Value = 19
$ws.Range("BO56").Value = 591
$ws.Range("BP56").Value = 104280177
$ws.Range("BQ56").Value = 3205
$ws.Range("BR56").Value = 6837
$ws.Range("BS56").Value = 10042
$ws.Range("BT56").Value = 0.033462074402905891
$ws.Range("BU56").Value = 0.75422556163723686
$ws.Range("BV56").Value = 45793
$ws.Range("BW56").Value = 0.75422556163723686
$ws.Range("BX56").Value = 45796
$ws.Range("BY56").Value = 0.062414107401757755
$ws.Range("BZ56").Value = 45807
$ws.Range("CA56").Value = 0.058626568499199519
$ws.Range("CB56").Value = 6.666666666666667
$ws.Range("CC56").Value = 585
$ws.Range("CD56").Value = -80651610
$ws.Range("CE56").Value = 0.048736046234529609
$ws.Range("CF56").Value = 10296
$ws.Range("CG56").Value = 41102
$ws.Range("CH56").Value = 51398
$ws.Range("CI56").Value = 0.37880206651598319
$ws.Range("CJ56").Value = 45793
$ws.Range("CK56").Value = 0.37880206651598319
$ws.Range("CL56").Value = 45807
$ws.Range("CM56").Value = 0.19299967710687763
$ws.Range("CN56").Value = 45814
$ws.Range("CO56").Value = 0.072069744914433326
$ws.Range("CP56").Value = 12.666666666666666
$ws.Range("CQ56").Value = 575
$ws.Range("CR56").Value = -78182175
$ws.Range("CS56").Value = 0.04724381937962658
$ws.Range("CT56").Value = -5866
$ws.Range("CU56").Value = 12623
$ws.Range("CV56").Value = 18489
$ws.Range("CW56").Value = 0.10041691015070639
$ws.Range("CX56").Value = 45793
$ws.Range("CY56").Value = 0.10041691015070639
$ws.Range("CZ56").Value = 45807
$ws.Range("DA56").Value = 0.16797206849308341
$ws.Range("DB56").Value = 45828
$ws.Range("DC56").Value = 0.32078019652231132
$ws.Range("DD56").Value = 17.333333333333332
$ws.Range("DE56").Value = 550
$ws.Range("DF56").Value = -64457250
$ws.Range("DG56").Value = 0.038950140191257601
$ws.Range("DH56").Value = -850
$ws.Range("DI56").Value = 109591
$ws.Range("DJ56").Value = 110441
$ws.Range("DK56").Value = 0
$ws.Range("DL56").Value = 45828
$ws.Range("DM56").Value = 0.94546695678143267
$ws.Range("DN56").Value = 45835
$ws.Range("DO56").Value = 0.0039848116387217886
$ws.Range("DP56").Value = 45838
$ws.Range("DQ56").Value = 0.050548231579845559
$ws.Range("DR56").Value = 41.666666666666664
$ws.Range("DS56").Value = 586
$ws.Range("DT56").Value = -55197098
$ws.Range("DU56").Value = 0.033354428016252394
$ws.Range("DV56").Value = -3837
$ws.Range("DW56").Value = 1639
$ws.Range("DX56").Value = 5476
$ws.Range("DY56").Value = 0.71445148895292987
$ws.Range("DZ56").Value = 45793
$ws.Range("EA56").Value = 0.71445148895292987
$ws.Range("EB56").Value = 45796
$ws.Range("EC56").Value = 0.061256484149855905
$ws.Range("ED56").Value = 45797
$ws.Range("EE56").Value = 0.061487031700288186
$ws.Range("EF56").Value = 3.3333333333333335
$ws.Range("EG56").Value = 587
$ws.Range("EH56").Value = -44533929
$ws.Range("EI56").Value = 0.026910902618673811
$ws.Range("EJ56").Value = -3559
$ws.Range("EK56").Value = -4038
$ws.Range("EL56").Value = 7597
$ws.Range("EM56").Value = 0.62627629474720137
$ws.Range("EN56").Value = 45793
$ws.Range("EO56").Value = 0.62627629474720137
$ws.Range("EP56").Value = 45796
$ws.Range("EQ56").Value = 0.11781891991634888
$ws.Range("ER56").Value = 45800
$ws.Range("ES56").Value = 0.061615819904047235
$ws.Range("ET56").Value = 4.333333333333333
$ws.Range("EU56").Value = 590
$ws.Range("EV56").Value = 493389860
$ws.Range("EW56").Value = -12757
$ws.Range("EX56").Value = 23172
$ws.Range("EY56").Value = 35929
$ws.Range("EZ56").Value = 0.10340925907762952
$ws.Range("FA56").Value = 388420010
$ws.Range("FB56").Value = 0.12463863840773352
$ws.Range("FC56").Value = 0.53833815101338367
$ws.Range("FD56").Value = 45793
$ws.Range("FE56").Value = 0.53833815101338367
$ws.Range("FF56").Value = 45807
$ws.Range("FG56").Value = 0.071388752603142147
$ws.Range("FH56").Value = 45828
$ws.Range("FI56").Value = 0.10913830108804126
$ws.Range("FJ56").Value = 17.333333333333332
$ws.Range("FK56").Value = -104969850
$ws.Range("FL56").Value = 0.063431039539466577
$ws.Range("FM56").Value = 0.43357783211083945
$ws.Range("FN56").Value = 45793
$ws.Range("FO56").Value = 0.43357783211083945
$ws.Range("FP56").Value = 45796
$ws.Range("FQ56").Value = 0.063985611106427226
$ws.Range("FR56").Value = 45856
$ws.Range("FS56").Value = 0.10882724896720344
$ws.Range("FT56").Value = 23
$ws.Range("FU56").Value = 600
$ws.Range("FV56").Value = 403924800
$ws.Range("FW56").Value = -11028
$ws.Range("FX56").Value = 1097
$ws.Range("FY56").Value = 12125
$ws.Range("FZ56").Value = 0.08465833548155953
$ws.Range("GA56").Value = 378898200
$ws.Range("GB56").Value = 0.12158322055328998
$ws.Range("GC56").Value = 0.18414972676038049
$ws.Range("GD56").Value = 45793
$ws.Range("GE56").Value = 0.18414972676038049
$ws.Range("GF56").Value = 45800
$ws.Range("GG56").Value = 0.071891077867353284
$ws.Range("GH56").Value = 45828
$ws.Range("GI56").Value = 0.34061602826300047
$ws.Range("GJ56").Value = 15
$ws.Range("GK56").Value = -25026600
$ws.Range("GL56").Value = 0.015123040131413107
$ws.Range("GM56").Value = 0.0074800412361247636
$ws.Range("GN56").Value = 45884
$ws.Range("GO56").Value = 0.20606075136055238
$ws.Range("GP56").Value = 45919
$ws.Range("GQ56").Value = 0.17918534679101436
$ws.Range("GR56").Value = 46038
$ws.Range("GS56").Value = 0.18450768382441082
$ws.Range("GT56").Value = 155
$ws.Range("GU56").Value = 595
$ws.Range("GV56").Value = 308935305
$ws.Range("GW56").Value = 24847
$ws.Range("GX56").Value = 1187
$ws.Range("GY56").Value = 26034
$ws.Range("GZ56").Value = 0.064749549155654507
$ws.Range("HA56").Value = 289618035
$ws.Range("HB56").Value = 0.092934443672773995
$ws.Range("HC56").Value = 0.37098076437125194
$ws.Range("HD56").Value = 45793
$ws.Range("HE56").Value = 0.37098076437125194
$ws.Range("HF56").Value = 45807
$ws.Range("HG56").Value = 0.19810047395701721
$ws.Range("HH56").Value = 45828
$ws.Range("HI56").Value = 0.097663496681068218
$ws.Range("HJ56").Value = 17.333333333333332
$ws.Range("HK56").Value = -19317270
$ws.Range("HL56").Value = 0.011673013890793895
$ws.Range("HM56").Value = 0.07983736832378488
$ws.Range("HN56").Value = 45796
$ws.Range("HO56").Value = 0.10977638144520421
$ws.Range("HP56").Value = 45919
$ws.Range("HQ56").Value = 0.18357666481857943
$ws.Range("HR56").Value = 46038
$ws.Range("HS56").Value = 0.14402759810263044
$ws.Range("HT56").Value = 125.66666666666667
$ws.Range("HU56").Value = 590
$ws.Range("HV56").Value = 36409
$ws.Range("HW56").Value = 361052
$ws.Range("HX56").Value = 3116369169
$ws.Range("HY56").Value = -1654865674
$ws.Range("HZ56").Value = 1461503495
$ws.Range("IA56").Value = 1.8831553629772129
$ws.Range("IB56").Value = 4771234843
$ws.Range("IC56").Value = 0.26945807433608232
$ws.Range("ID56").Value = 45793
$ws.Range("IE56").Value = 0.26945807433608232
$ws.Range("IF56").Value = 45807
$ws.Range("IG56").Value = 0.10618727785812324
$ws.Range("IH56").Value = 45828
$ws.Range("II56").Value = 0.18088509901502661

# Restore the active cell selection on the frozen (bottom-right) pane
$null = $ws.Range("G62").Select()
